$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.094.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.436.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.84%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.21%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +2.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.437.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.04%  "
$ws.Range("E10").Value = "  +4.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.75"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.75%  "
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.350"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +13.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.871.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.959.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000142"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.431.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "339.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.66%  "
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.174"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.45%  "
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +15.98%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.04%  "
$ws.Range("E31").Value = "  +10.78%  "
$ws.Range("E32").Value = "  +7.41%  "
$ws.Range("E33").Value = "  +13.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "174.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.398"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "375.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +19.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.23%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  +14.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "145.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.01%  "
$ws.Range("E46").Value = "  +10.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.595"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.35%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0952"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0519"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.66%  "
$ws.Range("E50").Value = "  +5.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.10%  "
